$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated interest_rate (column B) and apr (column C) figures, plus the
# refreshed updated_date (column F) for the Bankrate mortgage rate table.
# These columns store plain text values (e.g. "6.66%"), not numeric
# percentages or dates, so force text formatting before assigning.
$ws.Range("B2:C8").NumberFormat = "@"
$ws.Range("F2:F8").NumberFormat = "@"

$ws.Range("B2").Value = "6.66%"
$ws.Range("C2").Value = "6.72%"

$ws.Range("B3").Value = "6.33%"
$ws.Range("C3").Value = "6.42%"

$ws.Range("B4").Value = "5.86%"
$ws.Range("C4").Value = "5.95%"

$ws.Range("B5").Value = "5.77%"
$ws.Range("C5").Value = "5.85%"

$ws.Range("B6").Value = "6.73%"
$ws.Range("C6").Value = "6.80%"

$ws.Range("B7").Value = "6.79%"
$ws.Range("C7").Value = "6.84%"

$ws.Range("B8").Value = "6.65%"
$ws.Range("C8").Value = "6.69%"

$ws.Range("F2:F8").Value = "2025-08-13"
